# Lesson 4.2 The Observer Template for List Data - typo fixes
# (Examples/0-4-4-other-nats and Slides/L4.2)
#
# 1) Slide 33 ("Example 5: remove-first-even"): the last `cond` clause was
#    missing its `else`, and the continuation line's indentation didn't
#    line up with the rest of the clauses. Also the entrance animation on
#    that content box is removed.
# 2) Slide 37 ("Watch this work:"): the final two `insert`-trace lines
#    needed a little extra alignment padding.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# Slide 33: shape id=20 "Content Placeholder 2" (the `remove-first-even`
# code listing).
# ---------------------------------------------------------------------
$s33 = $p.Slides.Item(33)
$sh33 = $s33.Shapes.Item(3)
$tr33 = $sh33.TextFrame.TextRange

# paragraph 7: "    [(cons (first lst)"  -> "    [else (cons (first lst)"
$para7 = $tr33.Paragraphs(7, 1)
$run2 = $para7.Characters(5, 1)      # the lone "[" run
$run2.Text = "[else "

# paragraph 8: "           (remove-first-even (rest lst)))]))"
#   11 leading spaces -> 16 leading spaces (line up under "(cons ...")
$para8 = $tr33.Paragraphs(8, 1)
$lead = $para8.Characters(1, 11)
$lead.Text = "                "

# Remove the click-triggered fade-in animation on this shape.
$seq33 = $s33.TimeLine.MainSequence
for ($i = $seq33.Count; $i -ge 1; $i--) {
    $seq33.Item($i).Delete()
}

# ---------------------------------------------------------------------
# Slide 37: shape id=3 "Content Placeholder 2" (the `insert` trace).
# ---------------------------------------------------------------------
$s37 = $p.Slides.Item(37)
$sh37 = $s37.Shapes.Item(2)
$tr37 = $sh37.TextFrame.TextRange

# paragraph 1: "(insert 27 ..." -> "  (insert 27 ..."
$para1 = $tr37.Paragraphs(1, 1)
$open1 = $para1.Characters(1, 1)
$open1.Text = "  ("

# paragraph 4 (final result line): add extra spacing between "27" and the
# following "(" so it lines up with the line above it.
$para4 = $tr37.Paragraphs(4, 1)
$mid = $para4.Characters(27, 4)      # "27 (" -> "27   ("
$mid.Text = "27   ("
